# The upstream codeforIATI/codelists source data fixed a column-ordering
# bug in the Spanish SectorGroup codelist: the "category" and "group"
# columns were mislabeled/out of order. The fix swaps column D
# (codeforiati:category-name) with column E (codeforiati:group-name),
# and column F (codeforiati:group-code) with column G
# (codeforiati:category-code) -- for the header row and every data row.
#
# We swap via Range.Copy (instead of re-assigning .Value/.Value2) so that
# the underlying shared-string cell type is preserved and numeric-looking
# text such as "110"/"111" is not coerced into a literal number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

# Use far-away helper columns as scratch space for the 3-step swap so we
# never clobber real data while the swap is in progress.
$tempCol1 = "AA"
$tempCol2 = "AB"

$rangeD = $ws.Range("D1:D$lastRow")
$rangeE = $ws.Range("E1:E$lastRow")
$rangeF = $ws.Range("F1:F$lastRow")
$rangeG = $ws.Range("G1:G$lastRow")

$temp1 = $ws.Range("${tempCol1}1:${tempCol1}$lastRow")
$temp2 = $ws.Range("${tempCol2}1:${tempCol2}$lastRow")

# Swap D <-> E
$rangeD.Copy($temp1)
$rangeE.Copy($rangeD)
$temp1.Copy($rangeE)

# Swap F <-> G
$rangeF.Copy($temp2)
$rangeG.Copy($rangeF)
$temp2.Copy($rangeG)

$temp1.Clear()
$temp2.Clear()
